$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glucose-log rows continuing the existing 5-minute series (rows 402-415,
# i.e. A402:B415). Column A holds a literal timestamp string; column B holds
# a literal numeric-looking string (e.g. "20.2"). Both source columns in the
# existing data are plain shared-string text cells (t="s") with the default
# "Normal" / General style (no quotePrefix, no custom number format) - so we
# can't just assign the value directly, since Excel's normal text-entry
# coercion would store numeric-looking column-B values as real numbers (or,
# if forced to text via NumberFormat "@"/quote-prefix, would stamp a new
# cell style that the source file never had).
#
# Instead we write a TEXT() formula that evaluates to the exact literal
# string we want, then copy/paste-special the *values* on top of it. That
# bakes the formula's text result in as a plain literal (t="s") cell while
# leaving the cell's style completely untouched.

function Set-LiteralText {
    param($cell, [string]$text)
    $escaped = $text.Replace("""", """""")
    $cell.Formula = "=""" + $escaped + """"
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

function Set-LiteralNumberText {
    param($cell, [string]$numText)
    $cell.Formula = "=TEXT(" + $numText + ",""0.0"")"
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$rows = @(
    @{ Row = 402; A = "2026/02/12 23:33"; B = "20.2" },
    @{ Row = 403; A = "2026/02/12 23:38"; B = "20.3" },
    @{ Row = 404; A = "2026/02/12 23:43"; B = "19.4" },
    @{ Row = 405; A = "2026/02/12 23:48"; B = "19.5" },
    @{ Row = 406; A = "2026/02/12 23:53"; B = "19.3" },
    @{ Row = 407; A = "2026/02/12 23:58"; B = "18.0" },
    @{ Row = 408; A = "2026/02/13 00:03"; B = "18.6" },
    @{ Row = 409; A = "2026/02/13 00:08"; B = "18.5" },
    @{ Row = 410; A = "2026/02/13 00:13"; B = "17.9" },
    @{ Row = 411; A = "2026/02/13 00:18"; B = "17.2" },
    @{ Row = 412; A = "2026/02/13 00:23"; B = "17.7" },
    @{ Row = 413; A = "2026/02/13 00:28"; B = "16.7" },
    @{ Row = 414; A = "2026/02/13 00:33"; B = "16.5" },
    @{ Row = 415; A = "2026/02/13 00:38"; B = "16.3" }
)

foreach ($r in $rows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    Set-LiteralText $cellA $r.A

    $cellB = $ws.Cells.Item($r.Row, 2)
    Set-LiteralNumberText $cellB $r.B
}
